# Auto-generated edit script: update Leve profit-tracking values per scheduled-runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1214.1111
$ws.Range("J17").Value = 1214.1111
$ws.Range("L17").Value = 3642.3333
$ws.Range("N17").Value = -3978.3333

$ws.Range("H40").Value = 998.4231
$ws.Range("J40").Value = 998.4
$ws.Range("L40").Value = 998.4
$ws.Range("N40").Value = -1348.4

$ws.Range("H112").Value = 37038800
$ws.Range("I112").Value = 166667600
$ws.Range("J112").Value = 1997.619
$ws.Range("K112").Value = 500002800
$ws.Range("L112").Value = 5992.857
$ws.Range("M112").Value = -500001692
$ws.Range("N112").Value = -8208.857

$ws.Range("H121").Value = 814.6429000000001
$ws.Range("I121").Value = 300
$ws.Range("J121").Value = 854.2308
$ws.Range("K121").Value = 900
$ws.Range("L121").Value = 2562.6924
$ws.Range("M121").Value = 847
$ws.Range("N121").Value = -6056.6924

$ws.Range("H137").Value = 1451572.6
$ws.Range("I137").Value = 2274090
$ws.Range("K137").Value = 6822270
$ws.Range("M137").Value = -6819720

$ws.Range("H140").Value = 55077.3
$ws.Range("J140").Value = 55077.3
$ws.Range("L140").Value = 55077.3
$ws.Range("N140").Value = -65437.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24114.451
$ws.Range("I32").Value = 20891.627
$ws.Range("J32").Value = 37542.89
$ws.Range("K32").Value = 20891.627
$ws.Range("L32").Value = 37542.89
$ws.Range("M32").Value = -20604.627
$ws.Range("N32").Value = -38116.89

$ws.Range("H45").Value = 1100
$ws.Range("I45").Value = 1000
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 1000
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -623
$ws.Range("N45").Value = -1954

$ws.Range("H74").Value = 5006947.5
$ws.Range("I74").Value = 6309429
$ws.Range("K74").Value = 6309429
$ws.Range("M74").Value = -6308555

$ws.Range("H77").Value = 5006947.5
$ws.Range("I77").Value = 6309429
$ws.Range("K77").Value = 31547145
$ws.Range("M77").Value = -31542777

$ws.Range("H122").Value = 2593
$ws.Range("I122").Value = 2447.2856
$ws.Range("K122").Value = 7341.8568
$ws.Range("M122").Value = -4891.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 496.1
$ws.Range("I22").Value = 518
$ws.Range("J22").Value = 299
$ws.Range("K22").Value = 518
$ws.Range("L22").Value = 299
$ws.Range("M22").Value = -345
$ws.Range("N22").Value = -645

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 55937.145
$ws.Range("J98").Value = 55937.145
$ws.Range("L98").Value = 55937.145
$ws.Range("N98").Value = -60429.145

$ws.Range("H99").Value = 1622.4
$ws.Range("I99").Value = 1578
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 1578
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = -80
$ws.Range("N99").Value = -4796

$ws.Range("H126").Value = 1622.4
$ws.Range("I126").Value = 1578
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 4734
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -2264
$ws.Range("N126").Value = -10340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 480
$ws.Range("I32").Value = 75
$ws.Range("J32").Value = 530.625
$ws.Range("K32").Value = 225
$ws.Range("L32").Value = 1591.875
$ws.Range("M32").Value = 58
$ws.Range("N32").Value = -2157.875

$ws.Range("H33").Value = 264.25925
$ws.Range("I33").Value = 160.68182
$ws.Range("J33").Value = 720
$ws.Range("K33").Value = 964.0909199999999
$ws.Range("L33").Value = 4320
$ws.Range("M33").Value = -681.0909199999999
$ws.Range("N33").Value = -4886

$ws.Range("H41").Value = 1047.1177
$ws.Range("I41").Value = 260
$ws.Range("J41").Value = 1375.0834
$ws.Range("K41").Value = 780
$ws.Range("L41").Value = 4125.2502
$ws.Range("M41").Value = -442
$ws.Range("N41").Value = -4801.2502

$ws.Range("H44").Value = 590.9091
$ws.Range("I44").Value = 320
$ws.Range("J44").Value = 816.6667
$ws.Range("K44").Value = 960
$ws.Range("L44").Value = 2450.0001
$ws.Range("M44").Value = -562
$ws.Range("N44").Value = -3246.0001

$ws.Range("H46").Value = 1813.1177
$ws.Range("I46").Value = 274.33334
$ws.Range("J46").Value = 2142.8572
$ws.Range("K46").Value = 823.0000200000001
$ws.Range("L46").Value = 6428.571599999999
$ws.Range("M46").Value = -732.0000200000001
$ws.Range("N46").Value = -6610.571599999999

$ws.Range("H107").Value = 1163.2858
$ws.Range("I107").Value = 1142.4445
$ws.Range("J107").Value = 1200.8
$ws.Range("K107").Value = 3427.3335
$ws.Range("L107").Value = 3602.4
$ws.Range("M107").Value = -1507.3335
$ws.Range("N107").Value = -7442.4

$ws.Range("H129").Value = 4632313.5
$ws.Range("I129").Value = 2358.4285
$ws.Range("J129").Value = 7578648.5
$ws.Range("K129").Value = 7075.2855
$ws.Range("L129").Value = 22735945.5
$ws.Range("M129").Value = -2075.2855
$ws.Range("N129").Value = -22745945.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1642.5555
$ws.Range("I102").Value = 1564.3334
$ws.Range("J102").Value = 1799
$ws.Range("K102").Value = 1564.3334
$ws.Range("L102").Value = 1799
$ws.Range("M102").Value = 57.66660000000002
$ws.Range("N102").Value = -5043

$ws.Range("H122").Value = 2826.923
$ws.Range("I122").Value = 2600
$ws.Range("J122").Value = 2968.75
$ws.Range("K122").Value = 7800
$ws.Range("L122").Value = 8906.25
$ws.Range("M122").Value = -5350
$ws.Range("N122").Value = -13806.25

$ws.Range("H126").Value = 1702.6666
$ws.Range("I126").Value = 1381
$ws.Range("J126").Value = 1960
$ws.Range("K126").Value = 4143
$ws.Range("L126").Value = 5880
$ws.Range("M126").Value = -1673
$ws.Range("N126").Value = -10820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1424.091
$ws.Range("I82").Value = 1173.125
$ws.Range("J82").Value = 2093.3333
$ws.Range("K82").Value = 1173.125
$ws.Range("L82").Value = 2093.3333
$ws.Range("M82").Value = -812.125
$ws.Range("N82").Value = -2815.3333

$ws.Range("H85").Value = 1424.091
$ws.Range("I85").Value = 1173.125
$ws.Range("J85").Value = 2093.3333
$ws.Range("K85").Value = 1173.125
$ws.Range("L85").Value = 2093.3333
$ws.Range("M85").Value = 74.875
$ws.Range("N85").Value = -4589.3333

$ws.Range("H100").Value = 1956.9166
$ws.Range("I100").Value = 1900.6
$ws.Range("K100").Value = 1900.6
$ws.Range("M100").Value = -1359.6

$ws.Range("H122").Value = 3558.8838
$ws.Range("I122").Value = 4359.2856
$ws.Range("J122").Value = 3172.4827
$ws.Range("K122").Value = 13077.8568
$ws.Range("L122").Value = 9517.4481
$ws.Range("M122").Value = -10627.8568
$ws.Range("N122").Value = -14417.4481

$ws.Range("H132").Value = 25563.977
$ws.Range("I132").Value = 2044.3611
$ws.Range("J132").Value = 146522
$ws.Range("K132").Value = 6133.0833
$ws.Range("L132").Value = 439566
$ws.Range("M132").Value = -3603.0833
$ws.Range("N132").Value = -444626

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1060002
$ws.Range("J2").Value = 90003
$ws.Range("L2").Value = 90003
$ws.Range("N2").Value = -90227

$ws.Range("H122").Value = 742.4545000000001
$ws.Range("I122").Value = 757.2778
$ws.Range("J122").Value = 675.75
$ws.Range("K122").Value = 2271.8334
$ws.Range("L122").Value = 2027.25
$ws.Range("M122").Value = 178.1666
$ws.Range("N122").Value = -6927.25

